$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 263.33334
$ws.Cells.Item(12, 9).Value = 263.33334
$ws.Cells.Item(12, 11).Value = 263.33334
$ws.Cells.Item(12, 13).Value = -93.33334000000002
$ws.Cells.Item(49, 8).Value = 247
$ws.Cells.Item(132, 8).Value = 1097.3903
$ws.Cells.Item(132, 9).Value = 874.85
$ws.Cells.Item(132, 11).Value = 2624.55
$ws.Cells.Item(132, 13).Value = -94.55000000000018
$ws.Cells.Item(138, 8).Value = 1945.4125
$ws.Cells.Item(138, 9).Value = 1234.5106
$ws.Cells.Item(138, 10).Value = 2957.9092
$ws.Cells.Item(138, 11).Value = 3703.5318
$ws.Cells.Item(138, 12).Value = 8873.7276
$ws.Cells.Item(138, 13).Value = 1436.4682
$ws.Cells.Item(138, 14).Value = -19153.7276

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1888.75
$ws.Cells.Item(2, 9).Value = 2001.4286
$ws.Cells.Item(2, 11).Value = 2001.4286
$ws.Cells.Item(2, 13).Value = -1888.4286
$ws.Cells.Item(23, 8).Value = 20000.666
$ws.Cells.Item(23, 10).Value = 9999.6
$ws.Cells.Item(23, 12).Value = 9999.6
$ws.Cells.Item(23, 14).Value = -10517.6
$ws.Cells.Item(32, 8).Value = 670986.1
$ws.Cells.Item(32, 9).Value = 731054.3
$ws.Cells.Item(32, 11).Value = 731054.3
$ws.Cells.Item(32, 13).Value = -730767.3
$ws.Cells.Item(37, 8).Value = 12025
$ws.Cells.Item(37, 10).Value = 12025
$ws.Cells.Item(37, 12).Value = 12025
$ws.Cells.Item(37, 14).Value = -12571
$ws.Cells.Item(52, 8).Value = 61450
$ws.Cells.Item(52, 10).Value = 61450
$ws.Cells.Item(52, 12).Value = 61450
$ws.Cells.Item(52, 14).Value = -62086
$ws.Cells.Item(61, 8).Value = 1764.1571
$ws.Cells.Item(61, 9).Value = 1526.1455
$ws.Cells.Item(61, 10).Value = 2636.8667
$ws.Cells.Item(61, 11).Value = 1526.1455
$ws.Cells.Item(61, 12).Value = 2636.8667
$ws.Cells.Item(61, 13).Value = -1314.1455
$ws.Cells.Item(61, 14).Value = -3060.8667
$ws.Cells.Item(88, 8).Value = 2150
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 13).ClearContents()
$ws.Cells.Item(91, 8).Value = 2150
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 13).ClearContents()
$ws.Cells.Item(112, 8).Value = 63309.6
$ws.Cells.Item(112, 10).Value = 63309.6
$ws.Cells.Item(112, 12).Value = 63309.6
$ws.Cells.Item(112, 14).Value = -66263.60000000001
$ws.Cells.Item(116, 8).Value = 1888.75
$ws.Cells.Item(116, 9).Value = 2001.4286
$ws.Cells.Item(116, 11).Value = 2001.4286
$ws.Cells.Item(116, 13).Value = 292.5714
$ws.Cells.Item(121, 8).Value = 59980
$ws.Cells.Item(121, 10).Value = 59980
$ws.Cells.Item(121, 12).Value = 59980
$ws.Cells.Item(121, 14).Value = -63474
$ws.Cells.Item(122, 8).Value = 2321.2
$ws.Cells.Item(122, 9).Value = 1861.6
$ws.Cells.Item(122, 11).Value = 5584.799999999999
$ws.Cells.Item(122, 13).Value = -3134.799999999999
$ws.Cells.Item(127, 8).Value = 47500
$ws.Cells.Item(127, 10).Value = 47500
$ws.Cells.Item(127, 12).Value = 47500
$ws.Cells.Item(127, 14).Value = -57420
$ws.Cells.Item(129, 8).Value = 49666
$ws.Cells.Item(129, 10).Value = 49666
$ws.Cells.Item(129, 12).Value = 49666
$ws.Cells.Item(129, 14).Value = -59666
$ws.Cells.Item(131, 8).Value = 39741.668
$ws.Cells.Item(131, 10).Value = 39741.668
$ws.Cells.Item(131, 12).Value = 39741.668
$ws.Cells.Item(131, 14).Value = -49821.668
$ws.Cells.Item(136, 8).Value = 1764.1571
$ws.Cells.Item(136, 9).Value = 1526.1455
$ws.Cells.Item(136, 10).Value = 2636.8667
$ws.Cells.Item(136, 11).Value = 4578.4365
$ws.Cells.Item(136, 12).Value = 7910.6001
$ws.Cells.Item(136, 13).Value = -2028.4365
$ws.Cells.Item(136, 14).Value = -13010.6001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1888.75
$ws.Cells.Item(3, 9).Value = 2001.4286
$ws.Cells.Item(3, 11).Value = 2001.4286
$ws.Cells.Item(3, 13).Value = -1887.4286
$ws.Cells.Item(62, 8).Value = 30000
$ws.Cells.Item(62, 10).Value = 30000
$ws.Cells.Item(62, 12).Value = 30000
$ws.Cells.Item(62, 14).Value = -31372
$ws.Cells.Item(65, 8).Value = 30000
$ws.Cells.Item(65, 10).Value = 30000
$ws.Cells.Item(65, 12).Value = 90000
$ws.Cells.Item(65, 14).Value = -96864
$ws.Cells.Item(86, 8).Value = 1951.9166
$ws.Cells.Item(86, 10).Value = 2933.25
$ws.Cells.Item(86, 12).Value = 2933.25
$ws.Cells.Item(86, 14).Value = -5179.25
$ws.Cells.Item(89, 8).Value = 1951.9166
$ws.Cells.Item(89, 10).Value = 2933.25
$ws.Cells.Item(89, 12).Value = 14666.25
$ws.Cells.Item(89, 14).Value = -25898.25
$ws.Cells.Item(99, 8).Value = 2924.75
$ws.Cells.Item(99, 9).Value = 2899.6667
$ws.Cells.Item(99, 10).Value = 3000
$ws.Cells.Item(99, 11).Value = 2899.6667
$ws.Cells.Item(99, 12).Value = 3000
$ws.Cells.Item(99, 13).Value = -1401.6667
$ws.Cells.Item(99, 14).Value = -5996
$ws.Cells.Item(105, 8).Value = 125000000
$ws.Cells.Item(105, 9).Value = 125000000
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 125000000
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).ClearContents()
$ws.Cells.Item(105, 14).Value = -124998253
$ws.Cells.Item(110, 8).Value = 40702
$ws.Cells.Item(110, 10).Value = 40702
$ws.Cells.Item(110, 12).Value = 40702
$ws.Cells.Item(110, 14).Value = -48882

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1523.3667
$ws.Cells.Item(58, 9).Value = 1036.6471
$ws.Cells.Item(58, 10).Value = 2159.8462
$ws.Cells.Item(58, 11).Value = 1036.6471
$ws.Cells.Item(58, 12).Value = 2159.8462
$ws.Cells.Item(58, 13).Value = -833.6470999999999
$ws.Cells.Item(58, 14).Value = -2565.8462
$ws.Cells.Item(136, 8).Value = 1523.3667
$ws.Cells.Item(136, 9).Value = 1036.6471
$ws.Cells.Item(136, 10).Value = 2159.8462
$ws.Cells.Item(136, 11).Value = 3109.9413
$ws.Cells.Item(136, 12).Value = 6479.5386
$ws.Cells.Item(136, 13).Value = -559.9412999999995
$ws.Cells.Item(136, 14).Value = -11579.5386

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(114, 8).Value = 575
$ws.Cells.Item(114, 9).Value = 187.66667
$ws.Cells.Item(114, 10).Value = 907
$ws.Cells.Item(114, 11).Value = 563.00001
$ws.Cells.Item(114, 12).Value = 2721
$ws.Cells.Item(114, 13).Value = 2690.99999
$ws.Cells.Item(114, 14).Value = -9229

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 13).ClearContents()
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 12).ClearContents()
$ws.Cells.Item(98, 14).Value = 0
$ws.Cells.Item(107, 8).Value = 364.0909
$ws.Cells.Item(107, 9).Value = 411.2
$ws.Cells.Item(107, 10).Value = 324.83334
$ws.Cells.Item(107, 11).Value = 411.2
$ws.Cells.Item(107, 12).Value = 324.83334
$ws.Cells.Item(107, 13).Value = 1508.8
$ws.Cells.Item(107, 14).Value = -4164.83334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 12163.526
$ws.Cells.Item(22, 9).Value = 1733.3334
$ws.Cells.Item(22, 10).Value = 14119.1875
$ws.Cells.Item(22, 11).Value = 1733.3334
$ws.Cells.Item(22, 12).Value = 14119.1875
$ws.Cells.Item(22, 13).Value = -1438.3334
$ws.Cells.Item(22, 14).Value = -14709.1875
$ws.Cells.Item(27, 8).Value = 12163.526
$ws.Cells.Item(27, 9).Value = 1733.3334
$ws.Cells.Item(27, 10).Value = 14119.1875
$ws.Cells.Item(27, 11).Value = 1733.3334
$ws.Cells.Item(27, 12).Value = 14119.1875
$ws.Cells.Item(27, 13).Value = -1626.3334
$ws.Cells.Item(27, 14).Value = -14333.1875
$ws.Cells.Item(40, 8).Value = 102452.2
$ws.Cells.Item(40, 9).Value = 168768.67
$ws.Cells.Item(40, 10).Value = 2977.5
$ws.Cells.Item(40, 11).Value = 168768.67
$ws.Cells.Item(40, 12).Value = 2977.5
$ws.Cells.Item(40, 13).Value = -168632.67
$ws.Cells.Item(40, 14).Value = -3249.5
$ws.Cells.Item(46, 8).Value = 1749.75
$ws.Cells.Item(46, 9).Value = 2999
$ws.Cells.Item(46, 10).Value = 1333.3334
$ws.Cells.Item(46, 11).Value = 2999
$ws.Cells.Item(46, 12).Value = 1333.3334
$ws.Cells.Item(46, 13).Value = -2811
$ws.Cells.Item(46, 14).Value = -1709.3334
$ws.Cells.Item(55, 8).Value = 696.5
$ws.Cells.Item(55, 9).Value = 498
$ws.Cells.Item(55, 10).Value = 895
$ws.Cells.Item(55, 11).Value = 498
$ws.Cells.Item(55, 12).Value = 895
$ws.Cells.Item(55, 13).Value = -325
$ws.Cells.Item(55, 14).Value = -1241
$ws.Cells.Item(100, 8).Value = 3712
$ws.Cells.Item(100, 9).Value = 4980
$ws.Cells.Item(100, 11).Value = 4980
$ws.Cells.Item(100, 13).Value = -4439
$ws.Cells.Item(101, 8).Value = 17360
$ws.Cells.Item(101, 10).Value = 17360
$ws.Cells.Item(101, 12).Value = 17360
$ws.Cells.Item(101, 14).Value = -23850
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).ClearContents()
$ws.Cells.Item(110, 14).Value = 0

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 39931.168
$ws.Cells.Item(46, 10).Value = 39931.168
$ws.Cells.Item(46, 12).Value = 39931.168
$ws.Cells.Item(46, 14).Value = -40393.168
$ws.Cells.Item(134, 8).Value = 39931.168
$ws.Cells.Item(134, 10).Value = 39931.168
$ws.Cells.Item(134, 12).Value = 119793.504
$ws.Cells.Item(134, 14).Value = -124863.504
$ws.Cells.Item(136, 8).Value = 3002.5945
$ws.Cells.Item(136, 9).Value = 2400.3667
$ws.Cells.Item(136, 10).Value = 5583.5713
$ws.Cells.Item(136, 11).Value = 7201.1001
$ws.Cells.Item(136, 12).Value = 16750.7139
$ws.Cells.Item(136, 13).Value = -4651.1001
$ws.Cells.Item(136, 14).Value = -21850.7139
